# Registration Data Import XLS Template - unit test data update
# Replaces the placeholder "888-888-88XX" phone numbers in the
# "Individuals" sheet with realistic sample phone numbers, and
# updates the alternate "(541) 754-3010" phone_number_2 values to
# a different sample number as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Individuals")

# Data rows for individuals run from row 3 to row 29.
# Column H = phone_number_1, Column I = phone_number_2.
# Odd rows (3,5,7,...) use one pair of phone numbers, even rows
# (4,6,8,...) use another pair - this mirrors the original layout
# where the header/example rows alternated between two sample people.

For ($r = 3; $r -le 29; $r++) {
    If (($r % 2) -eq 1) {
        $ws.Cells.Item($r, 8).Value = "+44 1632 960852"
    } Else {
        $ws.Cells.Item($r, 8).Value = "+1-613-555-0182"
        $ws.Cells.Item($r, 9).Value = "+36 55 979 922"
    }
}
